$wb = $excel.ActiveWorkbook

# Work on the "Program" worksheet
$ws = $wb.Worksheets.Item("Program")

# Update the shared-string-backed cell contents
$ws.Range("A5").Value = "greenJavaGreatSelenium"
$ws.Range("B5").Value = "allAboutLogic"
$ws.Range("A6").Value = "KWAdfd"
$ws.Range("B7").Value = "BASicgraet"

# Make this sheet the active sheet and move the selection to B7
$ws.Activate()
$ws.Range("B7").Select()
